$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (Volume/Number and date range)
$ws.Range("A8").Value = "Volume 31   Number  11"
$ws.Range("C9").Value = "Report Covering the Week  3/11/2024  Through  3/17/2024"

# Weekly crime statistics table updates (rows 14-33)
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 2
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 100
$ws.Range("I14").Value = 16
$ws.Range("J14").Value = 13
$ws.Range("K14").Value = 23.076923076923
$ws.Range("L14").Value = 45.454545454545
$ws.Range("M14").Value = -33.333333333333
$ws.Range("N14").Value = -82.022471910112
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 3
$ws.Range("E15").Value = 133.333333333333
$ws.Range("F15").Value = 20
$ws.Range("G15").Value = 15
$ws.Range("H15").Value = 33.333333333333
$ws.Range("I15").Value = 39
$ws.Range("J15").Value = 54
$ws.Range("K15").Value = -27.777777777777
$ws.Range("L15").Value = -33.898305084745
$ws.Range("M15").Value = -11.363636363636
$ws.Range("N15").Value = -68.032786885245
$ws.Range("C16").Value = 38
$ws.Range("D16").Value = 32
$ws.Range("E16").Value = 18.75
$ws.Range("F16").Value = 177
$ws.Range("G16").Value = 162
$ws.Range("H16").Value = 9.259259259259
$ws.Range("I16").Value = 498
$ws.Range("J16").Value = 465
$ws.Range("K16").Value = 7.096774193548
$ws.Range("L16").Value = 1.425661914460
$ws.Range("M16").Value = -19.935691318328
$ws.Range("N16").Value = -84.981905910735
$ws.Range("C17").Value = 85
$ws.Range("D17").Value = 75
$ws.Range("E17").Value = 13.333333333333
$ws.Range("F17").Value = 274
$ws.Range("G17").Value = 265
$ws.Range("H17").Value = 3.396226415094
$ws.Range("I17").Value = 791
$ws.Range("J17").Value = 780
$ws.Range("K17").Value = 1.410256410256
$ws.Range("L17").Value = 5.046480743691
$ws.Range("M17").Value = 34.752981260647
$ws.Range("N17").Value = -46.554054054054
$ws.Range("C18").Value = 27
$ws.Range("D18").Value = 32
$ws.Range("E18").Value = -15.625
$ws.Range("F18").Value = 118
$ws.Range("G18").Value = 148
$ws.Range("H18").Value = -20.270270270270
$ws.Range("I18").Value = 383
$ws.Range("J18").Value = 456
$ws.Range("K18").Value = -16.008771929824
$ws.Range("L18").Value = -27.462121212121
$ws.Range("M18").Value = -29.595588235294
$ws.Range("N18").Value = -82.955051179350
$ws.Range("C19").Value = 92
$ws.Range("D19").Value = 108
$ws.Range("E19").Value = -14.814814814814
$ws.Range("F19").Value = 391
$ws.Range("G19").Value = 359
$ws.Range("H19").Value = 8.913649025069
$ws.Range("I19").Value = 1051
$ws.Range("J19").Value = 1165
$ws.Range("K19").Value = -9.785407725321
$ws.Range("L19").Value = -3.666361136571
$ws.Range("M19").Value = 43.383356070941
$ws.Range("N19").Value = -11.082910321489
$ws.Range("C20").Value = 28
$ws.Range("D20").Value = 40
$ws.Range("E20").Value = -30
$ws.Range("F20").Value = 102
$ws.Range("G20").Value = 123
$ws.Range("H20").Value = -17.073170731707
$ws.Range("I20").Value = 320
$ws.Range("J20").Value = 336
$ws.Range("K20").Value = -4.761904761904
$ws.Range("L20").Value = -14.666666666666
$ws.Range("M20").Value = 27.490039840637
$ws.Range("N20").Value = -84.826932195353
$ws.Range("C21").Value = 279
$ws.Range("D21").Value = 292
$ws.Range("E21").Value = -4.452054794520
$ws.Range("F21").Value = 1088
$ws.Range("G21").Value = 1075
$ws.Range("H21").Value = 1.209302325581
$ws.Range("I21").Value = 3098
$ws.Range("J21").Value = 3269
$ws.Range("K21").Value = -5.230957479351
$ws.Range("L21").Value = -6.348246674727
$ws.Range("M21").Value = 10.445632798574
$ws.Range("N21").Value = -70.621147463252
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 20
$ws.Range("G22").Value = 13
$ws.Range("H22").Value = 53.846153846153
$ws.Range("I22").Value = 70
$ws.Range("J22").Value = 57
$ws.Range("K22").Value = 22.807017543859
$ws.Range("L22").Value = -21.348314606741
$ws.Range("M22").Value = -17.647058823529
$ws.Range("C23").Value = 32
$ws.Range("D23").Value = 31
$ws.Range("E23").Value = 3.225806451612
$ws.Range("F23").Value = 102
$ws.Range("G23").Value = 113
$ws.Range("H23").Value = -9.734513274336
$ws.Range("I23").Value = 282
$ws.Range("J23").Value = 329
$ws.Range("K23").Value = -14.285714285714
$ws.Range("L23").Value = -10.191082802547
$ws.Range("M23").Value = 51.612903225806
$ws.Range("C24").Value = 237
$ws.Range("D24").Value = 203
$ws.Range("E24").Value = 16.748768472906
$ws.Range("F24").Value = 952
$ws.Range("G24").Value = 897
$ws.Range("H24").Value = 6.131549609810
$ws.Range("I24").Value = 2511
$ws.Range("J24").Value = 2493
$ws.Range("K24").Value = 0.722021660649
$ws.Range("L24").Value = 5.415617128463
$ws.Range("M24").Value = 33.492822966507
$ws.Range("C25").Value = 119
$ws.Range("D25").Value = 81
$ws.Range("E25").Value = 46.913580246913
$ws.Range("G25").Value = 353
$ws.Range("H25").Value = 31.728045325779
$ws.Range("I25").Value = 1117
$ws.Range("J25").Value = 997
$ws.Range("K25").Value = 12.036108324974
$ws.Range("L25").Value = 16.233090530697
$ws.Range("C26").Value = 137
$ws.Range("D26").Value = 106
$ws.Range("E26").Value = 29.245283018867
$ws.Range("F26").Value = 503
$ws.Range("G26").Value = 419
$ws.Range("H26").Value = 20.047732696897
$ws.Range("I26").Value = 1266
$ws.Range("J26").Value = 1183
$ws.Range("K26").Value = 7.016060862214
$ws.Range("L26").Value = 8.576329331046
$ws.Range("M26").Value = -12.326869806094
$ws.Range("D27").Value = 6
$ws.Range("E27").Value = 33.333333333333
$ws.Range("I27").Value = 71
$ws.Range("J27").Value = 70
$ws.Range("K27").Value = 1.428571428571
$ws.Range("L27").Value = -20.224719101123
$ws.Range("C28").Value = 10
$ws.Range("D28").Value = 8
$ws.Range("E28").Value = 25
$ws.Range("F28").Value = 40
$ws.Range("G28").Value = 38
$ws.Range("H28").Value = 5.263157894736
$ws.Range("I28").Value = 103
$ws.Range("J28").Value = 116
$ws.Range("K28").Value = -11.206896551724
$ws.Range("L28").Value = -5.504587155963
$ws.Range("C29").Value = 7
$ws.Range("D29").Value = 6
$ws.Range("E29").Value = 16.666666666666
$ws.Range("F29").Value = 15
$ws.Range("G29").Value = 15
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 37
$ws.Range("J29").Value = 41
$ws.Range("K29").Value = -9.756097560975
$ws.Range("L29").Value = -9.756097560975
$ws.Range("M29").Value = -45.588235294117
$ws.Range("N29").Value = -89.863013698630
$ws.Range("C30").Value = 5
$ws.Range("D30").Value = 6
$ws.Range("E30").Value = -16.666666666666
$ws.Range("F30").Value = 13
$ws.Range("G30").Value = 13
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 33
$ws.Range("J30").Value = 34
$ws.Range("K30").Value = -2.941176470588
$ws.Range("L30").Value = -13.157894736842
$ws.Range("M30").Value = -41.071428571428
$ws.Range("N30").Value = -90.350877192982
$ws.Range("D31").Value = 1
$ws.Range("F31").Value = 4
$ws.Range("G31").Value = 8
$ws.Range("H31").Value = -50
$ws.Range("I31").Value = 12
$ws.Range("J31").Value = 19
$ws.Range("K31").Value = -36.842105263157
$ws.Range("L31").Value = -33.333333333333
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 50
$ws.Range("J33").Value = 2
$ws.Range("K33").Value = 100
